$wb = $excel.ActiveWorkbook

# Clone the most recent weekly ranking sheet as a template, then rename it
# and place it at the end of the workbook (matching the existing date-sheet order).
$srcSheet = $wb.Worksheets.Item("2026-02-02")
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2026-02-09"

# rank, title, author, latest_episode for the week of 2026-02-09
$data = @(
    @(1, "宇崎ちゃんは遊びたい！", "丈(著者)", "第131話"),
    @(2, "ワンパンマン", "原作/ＯＮＥ 作画/村田雄介", "220撃目"),
    @(3, "転生したらスライムだった件", "原作：伏瀬 漫画：川上泰樹 キャラクター原案：みっつばー", "第132話　視察と研究成果"),
    @(4, "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～", "zunta(作画) はらわたさいぞう(原作)", "第36話：クソガキ①"),
    @(5, "魔王の俺が奴隷エルフを嫁にしたんだが、どう愛でればいい？", "原作／手島史詞 キャラクター原案／COMTA 漫画／板垣ハコ", "第76話"),
    @(6, "蜘蛛ですが、なにか？", "かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)", "第78話その1"),
    @(7, "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―", "光永康則", "第７４話『半球停止』⓵"),
    @(8, "勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが", "絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)", "第7話 後編"),
    @(9, "王子様の友達", "すけろく(著者)", "第34話"),
    @(10, "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。", "マツモトケンゴ", "第７３話　合評の戦いが始まった（１）"),
    @(11, "立ち飲み居酒屋でクーデレダウナー系女子大生が隣に来る話", "原作:剃り残し 漫画:大窪太一", "第1話"),
    @(12, "元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～", "沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)", "第83話その1"),
    @(13, "女友達は頼めば意外とヤらせてくれる", "ろくろ(漫画) 鏡遊(原作)", "第29話①"),
    @(14, "怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった", "菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)", "第19話前半"),
    @(15, "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～", "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり", "第５７話　和解する器用貧乏（１）"),
    @(16, "実は俺、最強でした？", "原作：澄守 彩 漫画：高橋 愛", "第138話　わたくしと戦いませんか？・前編"),
    @(17, "田舎で恋は難しい!?", "ねこうめ(著者)", "第4話"),
    @(18, "異世界魔王と召喚少女の奴隷魔術", "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大", "第133話　王都に戻ってみるⅡ（前編）"),
    @(19, "まんきつしたい常連さん", "しんみりん(著者)", "第53話前編"),
    @(20, "貞操逆転世界で頼めばヤれると噂の俺", "澄田佑貴(漫画) aaa168（スリーエー）(原作)", "休載イラスト"),
    @(21, "バキ外伝 烈海王は異世界転生しても一向にかまわんッッ", "板垣恵介 猪原賽 陸井栄史", "第89話　天罰"),
    @(22, "クセ強彼女は床にいざなう", "須河篤志(著者)", "休載イラストその２"),
    @(23, "異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～", "モリタ Ｕ４ nima", "第15話（１）　お嬢様と宰相と、とびきりの滋養食（１）"),
    @(24, "せっかく農家に転生したので勇者は目指しません", "マツオカヨシノリ(漫画) 月見里嘉助(原作) ゆーにっと(キャラクター原案)", "第2話④"),
    @(25, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第77話④"),
    @(26, "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！", "島知宏 音速炒飯 有都あらゆる", "第２７食　コカトリスのフライドチキン、パクパクですわ！（２）"),
    @(27, "ぽんドロイド！ はまさん", "はれやまはれぞう(著者)", "第16話"),
    @(28, "独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～", "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき", "第36話 独身貴族は生ハムの原木を買う（4）"),
    @(29, "聖者無双", "漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime", "第96話　戦闘準備・消えた魔物の死体①"),
    @(30, "よくわからないけれど異世界に転生していたようです", "内々けやき あし カオミン", "第147話 よくわからないけれど滾るみたいです（２）"),
    @(31, "悪友の俺がポンコツ騎士を見てられないんだが、どう世話を焼きゃいい？ ～まどめ外伝～", "原作：手島史詞 漫画：双葉もも キャラクター原案：COMTA 構成：板垣ハコ", "第19話"),
    @(32, "小林さんちのメイドラゴン", "クール教信者", "第159話"),
    @(33, "バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～", "板垣恵介 林たかあき", "第63話 同類の香り"),
    @(34, "追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。", "六志麻あさ 業務用餅 kisui", "第８０話"),
    @(35, "治癒魔法の間違った使い方 ~戦場を駆ける回復要員~", "九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)", "第85話その2"),
    @(36, "北斗の拳 世紀末ドラマ撮影伝", "原案/武論尊・原哲夫 漫画/倉尾宏", "第80話 映画北斗、観客の反応!!（後）"),
    @(37, "経験値貯蓄でのんびり傷心旅行 ～勇者と恋人に追放された戦士の無自覚ざまぁ～", "奏ヨシキ(著者) 徳川レモン(原作) riritto(キャラクターデザイン)", "第43話-2"),
    @(38, "落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～", "村上よしゆき 茨木野 あるてら", "第４５話　エリーゼ、死す（１）"),
    @(39, "アイドル辞めるけど結婚してくれますか!?", "三吉汐美(著者)", "第20話後半"),
    @(40, "ヴァージン・キラー忍法帖", "成田芋虫", "一之巻"),
    @(41, "淫獄団地", "搾精研究所(原作) 丈山雄為(漫画)", "コミックス告知"),
    @(42, "姫様“拷問”の時間です", "原作:春原ロビンソン　漫画:ひらけい", "拷問169"),
    @(43, "ふかふかダンジョン攻略記～俺の異世界転生冒険譚～", "KAKERU", "第73話「弓王2」(前半)"),
    @(44, "ハズレ枠の【状態異常スキル】で最強になった俺がすべてを蹂躙するまで", "鵜吉しょう（作画） 内々けやき（構成） 篠崎 芳（原作） KWKM（キャラクター原案）", "第61話　勇の剣"),
    @(45, "陰キャの俺が席替えでS級美少女に囲まれたら秘密の関係が始まった。", "星野 星野(原作) バラマツヒトミ(漫画) 黒兎 ゆう(キャラクター原案)", "第10話"),
    @(46, "ゲーム悪役貴族に転生した俺は、チート筋肉で無双する", "昼行燈（原作） しいたけ元帥（漫画）", "第38話"),
    @(47, "ライブダンジョン！", "ことりりょう(作画) dy冷凍(原作) Mika Pikazo(キャラクター原案)", "第93話後半"),
    @(48, "三枝さんはメガネ先輩と恋を描く", "セレビィ量産型(著者)", "第25話"),
    @(49, "魔王になったので、ダンジョン造って人外娘とほのぼのする", "遠野ノオト(作画) 流優(原作) だぶ竜(キャラクター原案)", "第86話後半"),
    @(50, "黄金の経験値", "原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)", "第22話（後編）")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
}

# restore the originally active sheet so the workbook view is unchanged
$wb.Worksheets.Item("Sheet1").Activate()
